$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031961604145692
$ws.Cells.Item(2, 4).Value = 1.0471271194131
$ws.Cells.Item(2, 5).Value = 1.042092597640028
$ws.Cells.Item(2, 6).Value = 1.05457607304981
$ws.Cells.Item(2, 9).Value = 1.038034499349782
$ws.Cells.Item(2, 10).Value = 1.037094446692479
$ws.Cells.Item(2, 11).Value = 1.049890626704791
$ws.Cells.Item(2, 12).Value = 1.044870261453213
$ws.Cells.Item(2, 13).Value = 1.057318903419241
$ws.Cells.Item(2, 14).Value = 1.038567239763178

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033008031795953
$ws.Cells.Item(3, 4).Value = 1.047707496914087
$ws.Cells.Item(3, 5).Value = 1.042964260252193
$ws.Cells.Item(3, 6).Value = 1.055394349435427
$ws.Cells.Item(3, 9).Value = 1.038136331595909
$ws.Cells.Item(3, 10).Value = 1.037782476176681
$ws.Cells.Item(3, 11).Value = 1.050283386389069
$ws.Cells.Item(3, 12).Value = 1.04555252553825
$ws.Cells.Item(3, 13).Value = 1.057950437564892
$ws.Cells.Item(3, 14).Value = 1.039256246328165

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033685856195599
$ws.Cells.Item(4, 4).Value = 1.048082961181144
$ws.Cells.Item(4, 5).Value = 1.043529096620848
$ws.Cells.Item(4, 6).Value = 1.055924257884599
$ws.Cells.Item(4, 9).Value = 1.03820088217984
$ws.Cells.Item(4, 10).Value = 1.038227826991
$ws.Cells.Item(4, 11).Value = 1.05053674850247
$ws.Cells.Item(4, 12).Value = 1.045994168939714
$ws.Cells.Item(4, 13).Value = 1.058358858853843
$ws.Cells.Item(4, 14).Value = 1.039702229591723

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.033970984441892
$ws.Cells.Item(5, 4).Value = 1.048240785732197
$ws.Cells.Item(5, 5).Value = 1.043766747102579
$ws.Cells.Item(5, 6).Value = 1.056147132493705
$ws.Cells.Item(5, 9).Value = 1.038227697396183
$ws.Cells.Item(5, 10).Value = 1.038415087576774
$ws.Cells.Item(5, 11).Value = 1.050643074007743
$ws.Cells.Item(5, 12).Value = 1.046179876258169
$ws.Cells.Item(5, 13).Value = 1.058530504347289
$ws.Cells.Item(5, 14).Value = 1.039889756109008

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.034018868735092
$ws.Cells.Item(6, 4).Value = 1.048267283907772
$ws.Cells.Item(6, 5).Value = 1.043806660962574
$ws.Cells.Item(6, 6).Value = 1.056184560003542
$ws.Cells.Item(6, 9).Value = 1.038232180896816
$ws.Cells.Item(6, 10).Value = 1.038446531490071
$ws.Cells.Item(6, 11).Value = 1.0506609154707
$ws.Cells.Item(6, 12).Value = 1.046211059677071
$ws.Cells.Item(6, 13).Value = 1.058559321128579
$ws.Cells.Item(6, 14).Value = 1.039921244676269

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.033689665422372
$ws.Cells.Item(7, 4).Value = 1.048085070123803
$ws.Cells.Item(7, 5).Value = 1.04353227136031
$ws.Cells.Item(7, 6).Value = 1.055927235550791
$ws.Cells.Item(7, 9).Value = 1.038201241751532
$ws.Cells.Item(7, 10).Value = 1.038230329039108
$ws.Cells.Item(7, 11).Value = 1.050538169969111
$ws.Cells.Item(7, 12).Value = 1.045996650212331
$ws.Cells.Item(7, 13).Value = 1.058361152606361
$ws.Cells.Item(7, 14).Value = 1.039704735193027

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.032315100780594
$ws.Cells.Item(8, 4).Value = 1.047323275959838
$ws.Cells.Item(8, 5).Value = 1.042387011165412
$ws.Cells.Item(8, 6).Value = 1.054852523623997
$ws.Cells.Item(8, 9).Value = 1.038069191451298
$ws.Cells.Item(8, 10).Value = 1.03732693802819
$ws.Cells.Item(8, 11).Value = 1.050023522348925
$ws.Cells.Item(8, 12).Value = 1.045100799411556
$ws.Cells.Item(8, 13).Value = 1.05753237884186
$ws.Cells.Item(8, 14).Value = 1.038800061263253

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029898445132247
$ws.Cells.Item(9, 4).Value = 1.045980371747345
$ws.Cells.Item(9, 5).Value = 1.040375191624094
$ws.Cells.Item(9, 6).Value = 1.05296209877839
$ws.Cells.Item(9, 9).Value = 1.037826254735398
$ws.Cells.Item(9, 10).Value = 1.035736227149875
$ws.Cells.Item(9, 11).Value = 1.049110732337714
$ws.Cells.Item(9, 12).Value = 1.043523562710194
$ws.Cells.Item(9, 13).Value = 1.056070319028603
$ws.Cells.Item(9, 14).Value = 1.037207091393034

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028291067852707
$ws.Cells.Item(10, 4).Value = 1.045084853661239
$ws.Cells.Item(10, 5).Value = 1.03903826887816
$ws.Cells.Item(10, 6).Value = 1.051704164838457
$ws.Cells.Item(10, 9).Value = 1.037657445578065
$ws.Cells.Item(10, 10).Value = 1.034676589382749
$ws.Cells.Item(10, 11).Value = 1.048498304084296
$ws.Cells.Item(10, 12).Value = 1.042473047082702
$ws.Cells.Item(10, 13).Value = 1.05509457689996
$ws.Cells.Item(10, 14).Value = 1.036145948818735

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027595944613255
$ws.Cells.Item(11, 4).Value = 1.044697047515558
$ws.Cells.Item(11, 5).Value = 1.038460398885297
$ws.Cells.Item(11, 6).Value = 1.051160042327465
$ws.Cells.Item(11, 9).Value = 1.037582732316279
$ws.Cells.Item(11, 10).Value = 1.03421796160643
$ws.Cells.Item(11, 11).Value = 1.04823220630134
$ws.Cells.Item(11, 12).Value = 1.042018405233466
$ws.Cells.Item(11, 13).Value = 1.054671840163739
$ws.Cells.Item(11, 14).Value = 1.035686669738369

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027337877686673
$ws.Cells.Item(12, 4).Value = 1.044552994671635
$ws.Cells.Item(12, 5).Value = 1.038245907597831
$ws.Cells.Item(12, 6).Value = 1.050958018368035
$ws.Cells.Item(12, 9).Value = 1.037554738053898
$ws.Cells.Item(12, 10).Value = 1.034047637894429
$ws.Cells.Item(12, 11).Value = 1.048133229928086
$ws.Cells.Item(12, 12).Value = 1.041849567719595
$ws.Cells.Item(12, 13).Value = 1.054514783029576
$ws.Cells.Item(12, 14).Value = 1.035516104147161

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027393227903908
$ws.Cells.Item(13, 4).Value = 1.044583894665229
$ws.Cells.Item(13, 5).Value = 1.038291909658278
$ws.Cells.Item(13, 6).Value = 1.05100134921459
$ws.Cells.Item(13, 9).Value = 1.03756075387869
$ws.Cells.Item(13, 10).Value = 1.034084171492456
$ws.Cells.Item(13, 11).Value = 1.048154466821638
$ws.Cells.Item(13, 12).Value = 1.04188578226366
$ws.Cells.Item(13, 13).Value = 1.054548473831259
$ws.Cells.Item(13, 14).Value = 1.035552689627088

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027574610009716
$ws.Cells.Item(14, 4).Value = 1.044685140141557
$ws.Cells.Item(14, 5).Value = 1.03844266579263
$ws.Cells.Item(14, 6).Value = 1.051143341172692
$ws.Cells.Item(14, 9).Value = 1.03758042324372
$ws.Cells.Item(14, 10).Value = 1.034203881957724
$ws.Cells.Item(14, 11).Value = 1.048224027646659
$ws.Cells.Item(14, 12).Value = 1.042004448317766
$ws.Cells.Item(14, 13).Value = 1.054658858456049
$ws.Cells.Item(14, 14).Value = 1.035672570094947

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027686383117942
$ws.Cells.Item(15, 4).Value = 1.044747520294275
$ws.Cells.Item(15, 5).Value = 1.038535572243555
$ws.Cells.Item(15, 6).Value = 1.051230838727991
$ws.Cells.Item(15, 9).Value = 1.037592510082949
$ws.Cells.Item(15, 10).Value = 1.034277643662487
$ws.Cells.Item(15, 11).Value = 1.048266868408278
$ws.Cells.Item(15, 12).Value = 1.042077567281871
$ws.Cells.Item(15, 13).Value = 1.054726865613757
$ws.Cells.Item(15, 14).Value = 1.035746436549791

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028337219497372
$ws.Cells.Item(16, 4).Value = 1.04511059036742
$ws.Cells.Item(16, 5).Value = 1.039076641929666
$ws.Cells.Item(16, 6).Value = 1.051740288624402
$ws.Cells.Item(16, 9).Value = 1.037662370019571
$ws.Cells.Item(16, 10).Value = 1.03470703131784
$ws.Cells.Item(16, 11).Value = 1.048515945003722
$ws.Cells.Item(16, 12).Value = 1.042503225263671
$ws.Cells.Item(16, 13).Value = 1.055122627723393
$ws.Cells.Item(16, 14).Value = 1.036176433984867

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028745708254569
$ws.Cells.Item(17, 4).Value = 1.045338324932936
$ws.Cells.Item(17, 5).Value = 1.039416316255242
$ws.Cells.Item(17, 6).Value = 1.052060007105224
$ws.Cells.Item(17, 9).Value = 1.037705758625855
$ws.Cells.Item(17, 10).Value = 1.034976429495499
$ws.Cells.Item(17, 11).Value = 1.048671940807975
$ws.Cells.Item(17, 12).Value = 1.0427702937197
$ws.Cells.Item(17, 13).Value = 1.055370817092995
$ws.Cells.Item(17, 14).Value = 1.036446214738838

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028984057887478
$ws.Cells.Item(18, 4).Value = 1.04547115458684
$ws.Cells.Item(18, 5).Value = 1.039614541495589
$ws.Cells.Item(18, 6).Value = 1.052246548492914
$ws.Cells.Item(18, 9).Value = 1.037730910331474
$ws.Cells.Item(18, 10).Value = 1.035133584375908
$ws.Cells.Item(18, 11).Value = 1.048762842412851
$ws.Cells.Item(18, 12).Value = 1.042926093168357
$ws.Cells.Item(18, 13).Value = 1.055515559095379
$ws.Cells.Item(18, 14).Value = 1.036603592797203

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029065343369298
$ws.Cells.Item(19, 4).Value = 1.0455164453227
$ws.Cells.Item(19, 5).Value = 1.039682147928958
$ws.Cells.Item(19, 6).Value = 1.052310163568307
$ws.Cells.Item(19, 9).Value = 1.037739459911053
$ws.Cells.Item(19, 10).Value = 1.035187173413176
$ws.Cells.Item(19, 11).Value = 1.048793822528045
$ws.Cells.Item(19, 12).Value = 1.042979220617157
$ws.Cells.Item(19, 13).Value = 1.055564908496177
$ws.Cells.Item(19, 14).Value = 1.036657257937051

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028701872477086
$ws.Cells.Item(20, 4).Value = 1.045313891576993
$ws.Cells.Item(20, 5).Value = 1.039379862171193
$ws.Cells.Item(20, 6).Value = 1.052025698637133
$ws.Cells.Item(20, 9).Value = 1.03770111958423
$ws.Cells.Item(20, 10).Value = 1.034947523600748
$ws.Cells.Item(20, 11).Value = 1.048655213038903
$ws.Cells.Item(20, 12).Value = 1.042741637422878
$ws.Cells.Item(20, 13).Value = 1.055344191056952
$ws.Cells.Item(20, 14).Value = 1.036417267794401

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027521193831948
$ws.Cells.Item(21, 4).Value = 1.044655325978814
$ws.Cells.Item(21, 5).Value = 1.038398267570645
$ws.Cells.Item(21, 6).Value = 1.051101525655824
$ws.Cells.Item(21, 9).Value = 1.037574637793226
$ws.Cells.Item(21, 10).Value = 1.034168629346554
$ws.Cells.Item(21, 11).Value = 1.048203547456366
$ws.Cells.Item(21, 12).Value = 1.041969503101823
$ws.Cells.Item(21, 13).Value = 1.054626353860662
$ws.Cells.Item(21, 14).Value = 1.035637267421025

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026779622783672
$ws.Cells.Item(22, 4).Value = 1.044241235212614
$ws.Cells.Item(22, 5).Value = 1.037781999509596
$ws.Cells.Item(22, 6).Value = 1.050520967675736
$ws.Cells.Item(22, 9).Value = 1.037493711352188
$ws.Cells.Item(22, 10).Value = 1.033679087673691
$ws.Cells.Item(22, 11).Value = 1.047918782250151
$ws.Cells.Item(22, 12).Value = 1.041484244459513
$ws.Cells.Item(22, 13).Value = 1.054174825682755
$ws.Cells.Item(22, 14).Value = 1.035147030542837

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027172670674119
$ws.Cells.Item(23, 4).Value = 1.044460754289078
$ws.Cells.Item(23, 5).Value = 1.038108609227801
$ws.Cells.Item(23, 6).Value = 1.050828683940067
$ws.Cells.Item(23, 9).Value = 1.037536744705403
$ws.Cells.Item(23, 10).Value = 1.033938585688026
$ws.Cells.Item(23, 11).Value = 1.04806981567664
$ws.Cells.Item(23, 12).Value = 1.041741468769676
$ws.Cells.Item(23, 13).Value = 1.054414207417807
$ws.Cells.Item(23, 14).Value = 1.035406897074116

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02872167971518
$ws.Cells.Item(24, 4).Value = 1.045324931971587
$ws.Cells.Item(24, 5).Value = 1.039396333896971
$ws.Cells.Item(24, 6).Value = 1.052041200987262
$ws.Cells.Item(24, 9).Value = 1.037703216249949
$ws.Cells.Item(24, 10).Value = 1.034960584870173
$ws.Cells.Item(24, 11).Value = 1.048662771870081
$ws.Cells.Item(24, 12).Value = 1.042754585899106
$ws.Cells.Item(24, 13).Value = 1.055356222286658
$ws.Cells.Item(24, 14).Value = 1.036430347612327

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030522553703771
$ws.Cells.Item(25, 4).Value = 1.046327596143136
$ws.Cells.Item(25, 5).Value = 1.040894544220168
$ws.Cells.Item(25, 6).Value = 1.053450411687119
$ws.Cells.Item(25, 9).Value = 1.037890270099264
$ws.Cells.Item(25, 10).Value = 1.036147319829794
$ws.Cells.Item(25, 11).Value = 1.049347403626826
$ws.Cells.Item(25, 12).Value = 1.043931148404863
$ws.Cells.Item(25, 13).Value = 1.056448484687271
$ws.Cells.Item(25, 14).Value = 1.037618767871711
